$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the very top; this pushes the existing
# row 1 (title) down to row 2, and the existing rows 3..22 down to 4..23
# (the previously-empty row 2 becomes the previously-empty row 3).
$ws.Rows(1).Insert()

$longStar = '****************************************************************************************************************************************'
$midStar  = '******************************************'
$star     = '*'

# New banner row at the top of the sheet.
$ws.Range("A1").Value = $longStar
$ws.Range("B1").Value = $longStar
$ws.Range("C1").Value = $midStar
$ws.Range("D1").Value = $star

# Mark the (now shifted) title row with an extra asterisk in column D.
$ws.Range("D2").Value = $star

# New banner row right below the title (fills the previously empty gap row).
$ws.Range("A3").Value = $longStar
$ws.Range("B3").Value = $longStar
$ws.Range("C3").Value = $midStar
$ws.Range("D3").Value = $star

# The row insert doesn't relocate existing conditional-formatting ranges,
# so shift them down by one row manually, preserving rule/dxf identity.
$ws.Range("A6:A7").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A7:A8"))
$ws.Range("A15").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A16"))

# Restore the selection to match the edited workbook.
$ws.Range("D4").Select()
